# "Generate Report for Archive"
#
# 1) The localization status string used throughout the workbook changes
#    from "Ready for handoff" -> "In Translation" (every cell that held the
#    old text — Overview!E2:F4 and the Status column (C) on the zh-cn /
#    de-de sheets — shares the same underlying string, so updating all of
#    them collapses back onto a single shared string, same as the source
#    edit did).
# 2) With the status text now shorter, the Status-ish columns (Overview
#    columns E/F, and column C on the zh-cn / de-de sheets) are narrower;
#    reflect that by tightening their column width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update every cell holding the old status text ---------------------
foreach ($addr in @("E2", "F2", "E3", "F3", "E4", "F4")) {
    $overview.Range($addr).Value = $newStatus
}

foreach ($sheet in @($zhcn, $dede)) {
    foreach ($addr in @("C2", "C3", "C4")) {
        $sheet.Range($addr).Value = $newStatus
    }
}

# --- Narrow the columns that showed the status text to their new fit ---
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
